# HighLevelSequenceDiagrams.pptx update
#   - Add the two global slide guides (horizontal @ 1488, vertical @ 2880)
#     that were dragged in while re-checking the diagram alignment.
#   - Rename the AddressBook-change event/handler to the new EzDo naming
#     used across the architecture + UI sequence diagrams.

$p = $ppt.ActivePresentation

# --- Slide guides -----------------------------------------------------
# (Best effort: PowerPoint exposes global guides via Presentation.Guides.
#  Guide.Position is in points; the new horizontal guide sits at 186pt
#  (1488 in the OOXML 1/8-pt units) and the vertical guide at the default
#  center, 360pt (2880 in OOXML units).)
try {
    $null = $p.Guides.Add(1, 186)
    $null = $p.Guides.Add(2, 360)
} catch {
    Write-Host "Guides.Add not available: $_"
}

# --- Rename AddressBookChangedEvent -> EzDoChangedEvent ----------------
# Walk every shape on the slide (including shapes nested in groups) and
# retarget just the substring that changed, leaving the surrounding run
# formatting (colors, sizes, the post(...)/handle...() wrapper text) intact.

$oldEvent = "AddressBookChangedEvent"
$oldHandler = "handleAddresssBookChangedEvent"
$newEvent = "EzDoChangedEvent"
$newHandler = "handleEzDoChangedEvent"

function Update-EzDoShape($sh) {
    if ($sh.Type -eq 6) {
        for ($j = 1; $j -le $sh.GroupItems.Count; $j++) {
            Update-EzDoShape($sh.GroupItems.Item($j))
        }
        return
    }

    if (-not $sh.HasTextFrame) { return }

    $tr = $sh.TextFrame.TextRange
    $full = $tr.Text

    if ($full -eq "post($oldEvent)") {
        $start = $full.IndexOf($oldEvent) + 1
        $tr.Characters($start, $oldEvent.Length).Text = $newEvent
    } elseif ($full -eq "$oldHandler()") {
        $start = $full.IndexOf($oldHandler) + 1
        $tr.Characters($start, $oldHandler.Length).Text = $newHandler
    }
}

for ($k = 1; $k -le $p.Slides.Count; $k++) {
    $s = $p.Slides.Item($k)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        Update-EzDoShape($s.Shapes.Item($i))
    }
}
